$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the date column (A2:A26) from 20250603 to 20250610
$ws.Range("A2:A26").Value = 20250610

# Update the selection to match the new selection range
$ws.Range("A2:D26").Select()
